$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (Bus index column skipped; columns follow sheet layout A..N)
$data = @{
    2 = @{ 2 = 1.02; 3 = 1.028925964812128; 4 = 1.038843469303303; 5 = 1.028830208089775; 6 = 1.047518908768287; 9 = 1.035813789077052; 10 = 1.034075577605533; 11 = 1.041630348385935; 12 = 1.031645884918717; 13 = 1.050281320580725; 14 = 1.015367270186937 }
    3 = @{ 2 = 1.02; 3 = 1.029710673985021; 4 = 1.03948803409868; 5 = 1.029492505515222; 6 = 1.048396060633527; 9 = 1.03598233251165; 10 = 1.034501946087754; 11 = 1.042085446852098; 12 = 1.032116593096212; 13 = 1.050970163590061; 14 = 1.015509660008074 }
    4 = @{ 2 = 1.02; 3 = 1.030219121120257; 4 = 1.03990572333368; 5 = 1.029922035725229; 6 = 1.048964803475374; 9 = 1.036090534481582; 10 = 1.034777846343683; 11 = 1.042379850819028; 12 = 1.032421468773213; 13 = 1.051416417675031; 14 = 1.015601766926854 }
    5 = @{ 2 = 1.02; 3 = 1.030433035159046; 4 = 1.040081464690328; 5 = 1.03010284280823; 6 = 1.049204180270428; 9 = 1.036135816711095; 10 = 1.034893836247104; 11 = 1.042503598853158; 12 = 1.032549707962351; 13 = 1.051604147179001; 14 = 1.015640481310949 }
    6 = @{ 2 = 1.02; 3 = 1.030468961765389; 4 = 1.040110980873929; 5 = 1.030133214696263; 6 = 1.049244388888983; 9 = 1.036143407702542; 10 = 1.034913311509298; 11 = 1.042524375514741; 12 = 1.032571243918289; 13 = 1.051635675025223; 14 = 1.01564698118098 }
    7 = @{ 2 = 1.02; 3 = 1.030221978812429; 4 = 1.039908071031682; 5 = 1.029924450766833; 6 = 1.048968000952794; 9 = 1.036091140354826; 10 = 1.034779396202412; 11 = 1.042381504424273; 12 = 1.032423182040399; 13 = 1.051418925639723; 14 = 1.015602284259872 }
    8 = @{ 2 = 1.02; 3 = 1.029191017445866; 4 = 1.039061174996206; 5 = 1.029053830648576; 6 = 1.047815103486557; 9 = 1.035870926104874; 10 = 1.034219667713963; 11 = 1.041784165767752; 12 = 1.031804900811811; 13 = 1.050514008193464; 14 = 1.01541539718036 }
    9 = @{ 2 = 1.02; 3 = 1.027379679639352; 4 = 1.037573611831102; 5 = 1.027527272866359; 6 = 1.045792579392228; 9 = 1.035476350159882; 10 = 1.033233497982454; 11 = 1.040731064815256; 12 = 1.030717739772449; 13 = 1.048923539973763; 14 = 1.015085878154152 }
    10 = @{ 2 = 1.02; 3 = 1.0261758245541; 4 = 1.036585226591829; 5 = 1.02651478229755; 6 = 1.044450423562777; 9 = 1.035208952928712; 10 = 1.032576224308934; 11 = 1.040028735522865; 12 = 1.029994616308769; 13 = 1.047866098998709; 14 = 1.014866092077952 }
    11 = @{ 2 = 1.02; 3 = 1.025655443213868; 4 = 1.03615805701446; 5 = 1.026077623753849; 6 = 1.04387075039505; 9 = 1.035092145322789; 10 = 1.032291674710407; 11 = 1.039724574725024; 12 = 1.029681905224921; 13 = 1.047408918930367; 14 = 1.01477090285362 }
    12 = @{ 2 = 1.02; 3 = 1.025462286848936; 4 = 1.035999510508889; 5 = 1.025915434366085; 6 = 1.043655659820562; 9 = 1.035048604846685; 10 = 1.032185989775638; 11 = 1.039611590141812; 12 = 1.029565812874547; 13 = 1.047239208718094; 14 = 1.014735542747969 }
    13 = @{ 2 = 1.02; 3 = 1.025503713341329; 4 = 1.036033513662052; 5 = 1.025950215876336; 6 = 1.043701787224699; 9 = 1.035057951335138; 10 = 1.032208659104034; 11 = 1.039635825954365; 12 = 1.029590712227631; 13 = 1.047275607271512; 14 = 1.014743127720206 }
    14 = @{ 2 = 1.02; 3 = 1.025639474054594; 4 = 1.036144948984459; 5 = 1.026064213225598; 6 = 1.043852966321531; 9 = 1.0350885493689; 10 = 1.03228293856028; 11 = 1.03971523549256; 12 = 1.02967230771228; 13 = 1.047394888441597; 14 = 1.01476798002633 }
    15 = @{ 2 = 1.02; 3 = 1.025723138844797; 4 = 1.036213624360202; 5 = 1.026134476091468; 6 = 1.043946142743368; 9 = 1.03510738158476; 10 = 1.03232870587686; 11 = 1.03976416161463; 12 = 1.029722589704555; 13 = 1.047468395716731; 14 = 1.014783292023463 }
    16 = @{ 2 = 1.02; 3 = 1.026210379581699; 4 = 1.03661359361386; 5 = 1.026543821713312; 6 = 1.044488926117711; 9 = 1.035216683577858; 10 = 1.032595110182965; 11 = 1.040048920792749; 12 = 1.030015378594429; 13 = 1.047896455414054; 14 = 1.014872409084038 }
    17 = @{ 2 = 1.02; 3 = 1.026516253969276; 4 = 1.036864701387914; 5 = 1.02680093121882; 6 = 1.044829799975853; 9 = 1.03528497250869; 10 = 1.032762234009925; 11 = 1.040227530971736; 12 = 1.03019914709136; 13 = 1.048165154344613; 14 = 1.014928304735331 }
    18 = @{ 2 = 1.02; 3 = 1.026694751637976; 4 = 1.037011246110678; 5 = 1.026951019950279; 6 = 1.045028769470612; 9 = 1.035324705547723; 10 = 1.032859719656374; 11 = 1.040331706580781; 12 = 1.030306375236676; 13 = 1.048321949191326; 14 = 1.014960905722573 }
    19 = @{ 2 = 1.02; 3 = 1.026755629314438; 4 = 1.037061227212396; 5 = 1.027002216810699; 6 = 1.045096637208958; 9 = 1.035338236716917; 10 = 1.032892960575595; 11 = 1.040367226932962; 12 = 1.030342943849427; 13 = 1.048375423542915; 14 = 1.014972021462936 }
    20 = @{ 2 = 1.02; 3 = 1.02648342760117; 4 = 1.036837751838336; 5 = 1.026773333261235; 6 = 1.044793212574526; 9 = 1.035277655960396; 10 = 1.03274430266327; 11 = 1.040208368264331; 12 = 1.030179426417067; 13 = 1.048136318517304; 14 = 1.014922307864425 }
    21 = @{ 2 = 1.02; 3 = 1.025599492127206; 4 = 1.036112130640348; 5 = 1.026030638565784; 6 = 1.043808441593959; 9 = 1.035079543226271; 10 = 1.032261064840265; 11 = 1.039691851509226; 12 = 1.029648278139655; 13 = 1.047359760146556; 14 = 1.014760661709265 }
    22 = @{ 2 = 1.02; 3 = 1.025044517094998; 4 = 1.035656618099681; 5 = 1.025564781587093; 6 = 1.043190584741319; 9 = 1.034954097446572; 10 = 1.031957289395722; 11 = 1.039367064536761; 12 = 1.029314686283415; 13 = 1.046872126015633; 14 = 1.014659013609341 }
    23 = @{ 2 = 1.02; 3 = 1.025338644334394; 4 = 1.035898025561081; 5 = 1.025811635818158; 6 = 1.043517997681559; 9 = 1.035020682215169; 10 = 1.032118320823233; 11 = 1.039539242849185; 12 = 1.029491494826847; 13 = 1.047130570816514; 14 = 1.014712900429375 }
    24 = @{ 2 = 1.02; 3 = 1.026498260157467; 4 = 1.03684992893961; 5 = 1.026785803216395; 6 = 1.044809744403059; 9 = 1.035280962298815; 10 = 1.03275240505125; 11 = 1.040217027081626; 12 = 1.030188337220324; 13 = 1.048149347978711; 14 = 1.014925017598499 }
    25 = @{ 2 = 1.02; 3 = 1.027847309027563; 4 = 1.037957605001606; 5 = 1.027921013879582; 6 = 1.046314367722711; 9 = 1.035579126846466; 10 = 1.033488421631137; 11 = 1.041003368948799; 12 = 1.030998512208949; 13 = 1.049334215871869; 14 = 1.015171087567408 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item([int]$r, [int]$c).Value = $data[$r][$c]
    }
}